$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to remain plain text so numeric-looking strings
    # (e.g. "558.65", "133.00") are not auto-coerced into numbers by Excel,
    # then restore the default "Normal" style so no stray formatting is left behind.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '69.203.22'
$ws.Range("E2").Value = '  -0.45%  '
Set-TextValue "D3" '2.468.07'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue "D5" '558.65'
$ws.Range("E5").Value = '  -2.07%  '
Set-TextValue "D6" '163.22'
$ws.Range("E6").Value = '  -2.20%  '
$ws.Range("E7").Value = '  +0.00%  '
Set-TextValue "D8" '0.504'
$ws.Range("E8").Value = '  -1.11%  '
Set-TextValue "D9" '2.467.33'
$ws.Range("E9").Value = '  -1.07%  '
Set-TextValue "D10" '0.152'
$ws.Range("E10").Value = '  -4.46%  '
$ws.Range("E11").Value = '  -0.58%  '
Set-TextValue "D12" '0.335'
$ws.Range("E12").Value = '  -4.44%  '
Set-TextValue "D13" '4.82'
$ws.Range("E13").Value = '  -1.26%  '
Set-TextValue "D14" '2.923.32'
$ws.Range("E14").Value = '  -0.94%  '
Set-TextValue "D15" '69.128.53'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  -3.10%  '
Set-TextValue "D17" '23.63'
$ws.Range("E17").Value = '  -2.57%  '
Set-TextValue "D18" '2.456.03'
$ws.Range("E18").Value = '  -1.66%  '
Set-TextValue "D19" '10.77'
$ws.Range("E19").Value = '  -4.11%  '
Set-TextValue "D20" '342.95'
$ws.Range("E20").Value = '  -2.87%  '
Set-TextValue "D21" '7.06'
$ws.Range("E21").Value = '  -5.49%  '
Set-TextValue "D22" '3.80'
$ws.Range("E22").Value = '  -2.88%  '
$ws.Range("E23").Value = '  -0.60%  '
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("E25").Value = '  +0.59%  '
Set-TextValue "D26" '67.12'
$ws.Range("E26").Value = '  -3.24%  '
Set-TextValue "D27" '3.68'
$ws.Range("E27").Value = '  -3.20%  '
Set-TextValue "D28" '2.596.66'
$ws.Range("E28").Value = '  -0.98%  '
Set-TextValue "D29" '0.998'
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("E30").Value = '  -5.12%  '
Set-TextValue "D31" '0.0₃0819'
$ws.Range("E31").Value = '  -6.34%  '
Set-TextValue "D32" '7.18'
$ws.Range("E32").Value = '  -5.31%  '
Set-TextValue "D33" '439.17'
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("E34").Value = '  -0.01%  '
Set-TextValue "D35" '1.14'
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("E36").Value = '  -5.91%  '
Set-TextValue "D37" '156.78'
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -3.98%  '
Set-TextValue "D41" '17.90'
$ws.Range("E41").Value = '  -1.37%  '
Set-TextValue "D42" '0.302'
$ws.Range("E42").Value = '  -3.48%  '
Set-TextValue "D43" '4.46'
$ws.Range("E43").Value = '  -3.76%  '
Set-TextValue "D44" '37.47'
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("E47").Value = '  -5.01%  '
Set-TextValue "D48" '133.00'
$ws.Range("E48").Value = '  -4.49%  '
Set-TextValue "D49" '3.36'
$ws.Range("E49").Value = '  -2.35%  '
Set-TextValue "D50" '0.0719'
$ws.Range("E50").Value = '  -0.60%  '
Set-TextValue "D51" '0.483'
$ws.Range("E51").Value = '  -4.38%  '
